$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# A new handoff was generated for 698895d0-0a15-4d92-beee-f20ac59001a9.md,
# so it now moves to "Ready for handoff" with a fresh handoff datetime, while
# c08e8002-d126-404b-9f55-d44a58821493.md remains "Handed back: in sync with
# en-US". The two entries also swap their row position (the 698895d0 entry
# moves down to row 3, c08e8002 moves up to row 2) on every sheet.
# ---------------------------------------------------------------------------

$missing = [System.Reflection.Missing]::Value

# ============================= Overview sheet =============================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value2 = "c08e8002-d126-404b-9f55-d44a58821493.md"
$ws1.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws1.Range("C2").Value2 = "Handed back: in sync with en-US"

$ws1.Range("A3").Value2 = "698895d0-0a15-4d92-beee-f20ac59001a9.md"
$ws1.Range("B3").Value2 = "Ready for handoff"
$ws1.Range("C3").Value2 = "Ready for handoff"

$hls1 = $ws1.Hyperlinks
$hls1.Delete()
$hls1.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f37cbf78070c8bff7107224bd2759a8738fc0929/e2e/c08e8002-d126-404b-9f55-d44a58821493.md", $missing, $missing, "c08e8002-d126-404b-9f55-d44a58821493.md")
$hls1.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f37cbf78070c8bff7107224bd2759a8738fc0929/e2e/698895d0-0a15-4d92-beee-f20ac59001a9.md", $missing, $missing, "698895d0-0a15-4d92-beee-f20ac59001a9.md")
$hls1.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f37cbf78070c8bff7107224bd2759a8738fc0929/.localization-config", $missing, $missing, ".localization-config")

# ============================== zh-cn sheet =================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value2 = "c08e8002-d126-404b-9f55-d44a58821493.md"
$ws2.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws2.Range("C2").Value2 = "c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.zh-cn.xlf"
$ws2.Range("D2").Value2 = "2016-03-09 16:37:35"
$ws2.Range("E2").Value2 = "c08e8002-d126-404b-9f55-d44a58821493.md"
$ws2.Range("F2").Value2 = "c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.zh-cn.xlf"
$ws2.Range("G2").Value2 = "2016-03-09 16:38:27"
$ws2.Range("H2").Value2 = "Include"

$ws2.Range("A3").Value2 = "698895d0-0a15-4d92-beee-f20ac59001a9.md"
$ws2.Range("B3").Value2 = "Ready for handoff"
$ws2.Range("C3").Value2 = "698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.zh-cn.xlf"
$ws2.Range("D3").Value2 = "2016-03-09 16:39:03"
$ws2.Range("E3").Value2 = "698895d0-0a15-4d92-beee-f20ac59001a9.md"
$ws2.Range("F3").Value2 = "698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.zh-cn.xlf"
$ws2.Range("G3").Value2 = "2016-03-09 16:38:27"
$ws2.Range("H3").Value2 = "Include"

$hls2 = $ws2.Hyperlinks
$hls2.Delete()
$hls2.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f37cbf78070c8bff7107224bd2759a8738fc0929/e2e/c08e8002-d126-404b-9f55-d44a58821493.md", $missing, $missing, "c08e8002-d126-404b-9f55-d44a58821493.md")
$hls2.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3197c7f8ed8988f4918e6e051d1f424b9e33ebc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.zh-cn.xlf", $missing, $missing, "c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.zh-cn.xlf")
$hls2.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/84255227afc8f892c8c32a85f8e43fa3fc53840d/e2e/c08e8002-d126-404b-9f55-d44a58821493.md", $missing, $missing, "c08e8002-d126-404b-9f55-d44a58821493.md")
$hls2.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c018864aa7fc0d9cef5c64879b899e329af6ccd8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.zh-cn.xlf", $missing, $missing, "c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.zh-cn.xlf")
$hls2.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f37cbf78070c8bff7107224bd2759a8738fc0929/e2e/698895d0-0a15-4d92-beee-f20ac59001a9.md", $missing, $missing, "698895d0-0a15-4d92-beee-f20ac59001a9.md")
$hls2.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3197c7f8ed8988f4918e6e051d1f424b9e33ebc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.zh-cn.xlf", $missing, $missing, "698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.zh-cn.xlf")
$hls2.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/84255227afc8f892c8c32a85f8e43fa3fc53840d/e2e/698895d0-0a15-4d92-beee-f20ac59001a9.md", $missing, $missing, "698895d0-0a15-4d92-beee-f20ac59001a9.md")
$hls2.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c018864aa7fc0d9cef5c64879b899e329af6ccd8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.zh-cn.xlf", $missing, $missing, "698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.zh-cn.xlf")
$hls2.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f37cbf78070c8bff7107224bd2759a8738fc0929/.localization-config", $missing, $missing, ".localization-config")

# ============================== de-de sheet ==================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value2 = "c08e8002-d126-404b-9f55-d44a58821493.md"
$ws3.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws3.Range("C2").Value2 = "c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.de-de.xlf"
$ws3.Range("D2").Value2 = "2016-03-09 16:37:48"
$ws3.Range("E2").Value2 = "c08e8002-d126-404b-9f55-d44a58821493.md"
$ws3.Range("F2").Value2 = "c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.de-de.xlf"
$ws3.Range("G2").Value2 = "2016-03-09 16:38:36"
$ws3.Range("H2").Value2 = "Include"

$ws3.Range("A3").Value2 = "698895d0-0a15-4d92-beee-f20ac59001a9.md"
$ws3.Range("B3").Value2 = "Ready for handoff"
$ws3.Range("C3").Value2 = "698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.de-de.xlf"
$ws3.Range("D3").Value2 = "2016-03-09 16:39:07"
$ws3.Range("E3").Value2 = "698895d0-0a15-4d92-beee-f20ac59001a9.md"
$ws3.Range("F3").Value2 = "698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.de-de.xlf"
$ws3.Range("G3").Value2 = "2016-03-09 16:38:36"
$ws3.Range("H3").Value2 = "Include"

$hls3 = $ws3.Hyperlinks
$hls3.Delete()
$hls3.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f37cbf78070c8bff7107224bd2759a8738fc0929/e2e/c08e8002-d126-404b-9f55-d44a58821493.md", $missing, $missing, "c08e8002-d126-404b-9f55-d44a58821493.md")
$hls3.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7e7a1b78554c626a2c5f7c79cbed54211d2c035f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.de-de.xlf", $missing, $missing, "c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.de-de.xlf")
$hls3.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/03c1b90b36db3d5f88d5ef7d84715b48dc2585ce/e2e/c08e8002-d126-404b-9f55-d44a58821493.md", $missing, $missing, "c08e8002-d126-404b-9f55-d44a58821493.md")
$hls3.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b24f1960b57d63670f11c12be78f72cbade0ff9f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.de-de.xlf", $missing, $missing, "c08e8002-d126-404b-9f55-d44a58821493.768a826f209fe896d36a25b0f8984b9a2ec4dd65.de-de.xlf")
$hls3.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f37cbf78070c8bff7107224bd2759a8738fc0929/e2e/698895d0-0a15-4d92-beee-f20ac59001a9.md", $missing, $missing, "698895d0-0a15-4d92-beee-f20ac59001a9.md")
$hls3.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7e7a1b78554c626a2c5f7c79cbed54211d2c035f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.de-de.xlf", $missing, $missing, "698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.de-de.xlf")
$hls3.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/03c1b90b36db3d5f88d5ef7d84715b48dc2585ce/e2e/698895d0-0a15-4d92-beee-f20ac59001a9.md", $missing, $missing, "698895d0-0a15-4d92-beee-f20ac59001a9.md")
$hls3.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b24f1960b57d63670f11c12be78f72cbade0ff9f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.de-de.xlf", $missing, $missing, "698895d0-0a15-4d92-beee-f20ac59001a9.8206f6eaccf7882b1a486c35d54f653b24f4b131.de-de.xlf")
$hls3.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f37cbf78070c8bff7107224bd2759a8738fc0929/.localization-config", $missing, $missing, ".localization-config")
